$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New order rows appended to the order history sheet (person/order details added)
$ws.Range("A17:A26").NumberFormat = "@"  # keep numeric-looking Order IDs as text, matching column A

$orders = @(
    @{ Row = 17; OrderId = '20250308120156'; OrderDate = '2025-03-08 12:01:56'; Items = '[{''Item Name'': ''Dal Rice'', ''Price'': 120.0}, {''Item Name'': ''South Indian Thali'', ''Price'': 200.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'; Total = 420; Payment = 'Paid'; OrderType = 'Delivery'; Delivery = 'Delivered' }
    @{ Row = 18; OrderId = '20250308120254'; OrderDate = '2025-03-08 12:02:54'; Items = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'; Total = 300; Payment = 'Paid'; OrderType = 'Take Away'; Delivery = 'Delivered' }
    @{ Row = 19; OrderId = '20250308121654'; OrderDate = '2025-03-08 12:16:54'; Items = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'; Total = 600; Payment = 'Not Paid'; OrderType = 'Take Away'; Delivery = 'Pending' }
    @{ Row = 20; OrderId = '20250308121701'; OrderDate = '2025-03-08 12:17:01'; Items = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'; Total = 600; Payment = 'Not Paid'; OrderType = 'Take Away'; Delivery = 'Pending' }
    @{ Row = 21; OrderId = '20250308121955'; OrderDate = '2025-03-08 12:19:55'; Items = '[{''Item Name'': ''Veg Thali'', ''Price'': 150.0}, {''Item Name'': ''Veg Thali'', ''Price'': 150.0}, {''Item Name'': ''Veg Thali'', ''Price'': 150.0}]'; Total = 450; Payment = 'Not Paid'; OrderType = 'Take Away'; Delivery = 'Pending' }
    @{ Row = 22; OrderId = '20250308121955'; OrderDate = '2025-03-08 12:19:55'; Items = '[{''Item Name'': ''Veg Thali'', ''Price'': 150.0}, {''Item Name'': ''Veg Thali'', ''Price'': 150.0}, {''Item Name'': ''Veg Thali'', ''Price'': 150.0}]'; Total = 450; Payment = 'Not Paid'; OrderType = 'Take Away'; Delivery = 'Pending' }
    @{ Row = 23; OrderId = '20250308122353'; OrderDate = '2025-03-08 12:23:53'; Items = '[{''Item Name'': ''Dal Rice'', ''Price'': 120.0}, {''Item Name'': ''Dal Rice'', ''Price'': 120.0}, {''Item Name'': ''Dal Rice'', ''Price'': 120.0}, {''Item Name'': ''Dal Rice'', ''Price'': 120.0}]'; Total = 480; Payment = 'Not Paid'; OrderType = 'Take Away'; Delivery = 'Pending' }
    @{ Row = 24; OrderId = '20250308122402'; OrderDate = '2025-03-08 12:24:02'; Items = '[{''Item Name'': ''Dal Rice'', ''Price'': 120.0}, {''Item Name'': ''Dal Rice'', ''Price'': 120.0}, {''Item Name'': ''Dal Rice'', ''Price'': 120.0}, {''Item Name'': ''Dal Rice'', ''Price'': 120.0}]'; Total = 480; Payment = 'Not Paid'; OrderType = 'Take Away'; Delivery = 'Pending' }
    @{ Row = 25; OrderId = '20250308122537'; OrderDate = '2025-03-08 12:25:37'; Items = '[{''Item Name'': ''Paneer Special'', ''Price'': 180.0}, {''Item Name'': ''Paneer Special'', ''Price'': 180.0}]'; Total = 360; Payment = 'Not Paid'; OrderType = 'Take Away'; Delivery = 'Pending' }
    @{ Row = 26; OrderId = '20250308122656'; OrderDate = '2025-03-08 12:26:56'; Items = '[{''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}, {''Item Name'': ''Roti Sabzi'', ''Price'': 100.0}]'; Total = 300; Payment = 'Paid'; OrderType = 'Take Away'; Delivery = 'Delivered' }
)

foreach ($o in $orders) {
    $r = $o.Row
    $ws.Cells.Item($r, 1).Value2 = $o.OrderId
    $ws.Cells.Item($r, 2).Value2 = $o.OrderDate
    $ws.Cells.Item($r, 3).Value2 = $o.Items
    $ws.Cells.Item($r, 4).Value2 = $o.Total
    $ws.Cells.Item($r, 5).Value2 = $o.Payment
    $ws.Cells.Item($r, 6).Value2 = $o.OrderType
    $ws.Cells.Item($r, 7).Value2 = $o.Delivery
}
